$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three "Requisitos" entries need to be reordered so that the
# "LOM3246" entry moves to the top (row 23), pushing the other two
# down by one position (row 24, row 25).

$lom3246 = "LOM3246 -  Técnicas de Caracterização de Materiais  (Indicação de Conjunto)`n"
$lob1021 = "LOB1021 -  Física IV  (Requisito)`n"
$lom3016 = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)`n"

$ws.Range("B23").Value = $lom3246
$ws.Range("C23").Value = $lom3246

$ws.Range("B24").Value = $lob1021
$ws.Range("C24").Value = $lob1021

$ws.Range("B25").Value = $lom3016
$ws.Range("C25").Value = $lom3016
